# This script updates the cryptocurrency price/volume table (columns B-E, rows 2-51)
# on the active worksheet to reflect refreshed data from coinranking.com.
#
# Note: several "Price" values in column D look numeric but are really text
# (e.g. "306.73", "21.13") that must stay literal strings rather than be
# reinterpreted as numbers by Excel (which would introduce binary floating-point
# rounding noise, e.g. 306.73 -> 306.73000000000002). We force those through as
# text by prefixing with a leading apostrophe, Excel's standard 'treat as text'
# marker, which is stripped from the stored value automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.072.88'
$ws.Range("E2").Value = '  +0.39%  '

$ws.Range("D3").Value = '1.893.69'
$ws.Range("E3").Value = '  +1.14%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = '''306.73'
$ws.Range("E5").Value = '  +0.44%  '

$ws.Range("E6").Value = '  -0.04%  '

$ws.Range("D7").Value = '''0.5187'
$ws.Range("E7").Value = '  +2.53%  '

$ws.Range("D8").Value = '''0.3757'
$ws.Range("E8").Value = '  +2.66%  '

$ws.Range("D9").Value = '''0.07216'
$ws.Range("E9").Value = '  +0.27%  '

$ws.Range("D10").Value = '''21.13'
$ws.Range("E10").Value = '  +2.05%  '

$ws.Range("D11").Value = '''0.8969'
$ws.Range("E11").Value = '  +0.37%  '

$ws.Range("D12").Value = '''0.07659'
$ws.Range("E12").Value = '  +1.73%  '

$ws.Range("D13").Value = '1.898.37'
$ws.Range("E13").Value = '  +1.25%  '

$ws.Range("D14").Value = '''94.21'
$ws.Range("E14").Value = '  -0.58%  '

$ws.Range("E15").Value = '  -0.16%  '

$ws.Range("E16").Value = '  -0.09%  '

$ws.Range("D17").Value = '''0.000008514'
$ws.Range("E17").Value = '  -0.20%  '

$ws.Range("D18").Value = '''14.45'
$ws.Range("E18").Value = '  +1.54%  '

$ws.Range("D19").Value = '''0.9999'
$ws.Range("E19").Value = '  -0.04%  '

$ws.Range("D20").Value = '27.125.69'
$ws.Range("E20").Value = '  +0.40%  '

$ws.Range("E21").Value = '  +0.74%  '

$ws.Range("D22").Value = '2.129.23'
$ws.Range("E22").Value = '  +0.37%  '

$ws.Range("E23").Value = '  +1.92%  '

$ws.Range("D24").Value = '''6.414'
$ws.Range("E24").Value = '  -0.07%  '

$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = '''146.75'
$ws.Range("E25").Value = '  -0.98%  '

$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").Value = '''2.281'
$ws.Range("E26").Value = '  +9.90%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '''18.05'
$ws.Range("E27").Value = '  +0.85%  '

$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = '''1.729'
$ws.Range("E28").Value = '  -3.51%  '

$ws.Range("D29").Value = '''114.43'
$ws.Range("E29").Value = '  +0.94%  '

$ws.Range("D30").Value = '''4.988'
$ws.Range("E30").Value = '  +6.40%  '

$ws.Range("D31").Value = '''4.786'
$ws.Range("E31").Value = '  +1.76%  '

$ws.Range("D32").Value = '''0.09189'
$ws.Range("E32").Value = '  +0.47%  '

$ws.Range("D33").Value = '''0.05045'
$ws.Range("E33").Value = '  -1.55%  '

$ws.Range("E34").Value = '  +6.71%  '

$ws.Range("D35").Value = '''0.7725'
$ws.Range("E35").Value = '  +2.87%  '

$ws.Range("D36").Value = '''2.975'
$ws.Range("E36").Value = '  -0.27%  '

$ws.Range("E37").Value = '  +1.72%  '

$ws.Range("D38").Value = '''2.599'
$ws.Range("E38").Value = '  +1.44%  '

$ws.Range("D39").Value = '''0.5600'
$ws.Range("E39").Value = '  -0.98%  '

$ws.Range("D40").Value = '''0.01990'
$ws.Range("E40").Value = '  -0.60%  '

$ws.Range("D41").Value = '''1.074'
$ws.Range("E41").Value = '  +0.07%  '

$ws.Range("D42").Value = '''119.64'
$ws.Range("E42").Value = '  +3.56%  '

$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = '''8.975'
$ws.Range("E43").Value = '  +5.42%  '

$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '''6.619'
$ws.Range("E44").Value = '  +0.19%  '

$ws.Range("D45").Value = '''0.1516'
$ws.Range("E45").Value = '  +2.91%  '

$ws.Range("D46").Value = '''0.4821'
$ws.Range("E46").Value = '  +1.88%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '''10.18'
$ws.Range("E47").Value = '  +0.82%  '

$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D48").Value = '''0.9999'
$ws.Range("E48").Value = '  +0.01%  '

$ws.Range("D49").Value = '''1.594'
$ws.Range("E49").Value = '  +1.91%  '

$ws.Range("D50").Value = '''37.40'
$ws.Range("E50").Value = '  +1.45%  '

$ws.Range("D51").Value = '''63.97'
$ws.Range("E51").Value = '  +1.26%  '
